$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; New="89÷5=17, 4"},
    @{Row=1;  Col=2; New="91÷2=45, 1"},
    @{Row=1;  Col=3; New="81÷2=40, 1"},
    @{Row=1;  Col=4; New="52÷6=8, 4"},
    @{Row=1;  Col=5; New="46÷6=7, 4"},

    @{Row=5;  Col=1; New="50÷9=5, 5"},
    @{Row=5;  Col=2; New="67÷4=16, 3"},
    @{Row=5;  Col=3; New="81÷5=16, 1"},
    @{Row=5;  Col=4; New="83÷4=20, 3"},
    @{Row=5;  Col=5; New="95÷6=15, 5"},

    @{Row=9;  Col=1; New="32÷7=4, 4"},
    @{Row=9;  Col=2; New="11÷4=2, 3"},
    @{Row=9;  Col=3; New="37÷3=12, 1"},
    @{Row=9;  Col=4; New="66÷6=11, 0"},
    @{Row=9;  Col=5; New="96÷3=32, 0"},

    @{Row=13; Col=1; New="24÷7=3, 3"},
    @{Row=13; Col=2; New="62÷6=10, 2"},
    @{Row=13; Col=3; New="15÷6=2, 3"},
    @{Row=13; Col=4; New="36÷7=5, 1"},
    @{Row=13; Col=5; New="61÷2=30, 1"},

    @{Row=17; Col=1; New="40÷5=8, 0"},
    @{Row=17; Col=2; New="64÷2=32, 0"},
    @{Row=17; Col=3; New="14÷3=4, 2"},
    @{Row=17; Col=4; New="62÷9=6, 8"},
    @{Row=17; Col=5; New="10÷3=3, 1"}
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    $cell.Range.Text = $item.New
}
